$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels for columns B and C (column A header stays "material name")
$ws.Range("B1").Value = "Conductivity"
$ws.Range("C1").Value = "permittivity"

# Move selection to C1, matching the final state in the workbook
$ws.Range("C1").Select()
